# repull data, push all data, mean calculation
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = -2
$ws.Range("F10").Value = -2
$ws.Range("F12").Value = 0
$ws.Range("F13").Value = -1
$ws.Range("F15").Value = 0
$ws.Range("F16").Value = -6
$ws.Range("F18").Value = 5
